# Weekly update: insert the newest week's record at row 126, pushing the
# existing historical rows (126-188) down by one (to 127-189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 126; everything from old row 126
# downward shifts to row+1.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with this week's data.
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value = "Maule"
$ws.Cells.Item(126, 4).Value = 44455
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 100112023
$ws.Cells.Item(126, 7).Value = "Brócoli"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 600
$ws.Cells.Item(126, 12).Value = 600
$ws.Cells.Item(126, 13).Value = 600
$ws.Cells.Item(126, 14).Value = "`$/unidad"
$ws.Cells.Item(126, 15).Value = "Región del Maule"
$ws.Cells.Item(126, 16).Value = 600
$ws.Cells.Item(126, 17).Value = 1
$ws.Cells.Item(126, 18).Value = "Hortaliza"
